$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.307.90"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "2.644.87"
$ws.Range("E3").Value = "  +2.57%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'594.12"
$ws.Range("E5").Value = "  +2.00%  "
$ws.Range("D6").Value = "'142.89"
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.586"
$ws.Range("E8").Value = "  -0.89%  "
$ws.Range("D9").Value = "2.642.66"
$ws.Range("E9").Value = "  +2.53%  "
$ws.Range("D10").Value = "'0.106"
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").Value = "'5.65"
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "'0.354"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").Value = "'27.43"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").Value = "3.118.68"
$ws.Range("E15").Value = "  +2.57%  "
$ws.Range("D16").Value = "63.284.69"
$ws.Range("E16").Value = "  +0.97%  "
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").Value = "2.658.39"
$ws.Range("E18").Value = "  +2.95%  "
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("D20").Value = "'338.74"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "'4.36"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "'6.71"
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'67.09"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").Value = "'1.67"
$ws.Range("E25").Value = "  +4.90%  "
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").Value = "'8.37"
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("D30").Value = "'7.83"
$ws.Range("E30").Value = "  -2.40%  "
$ws.Range("D31").Value = "'527.77"
$ws.Range("E31").Value = "  +16.20%  "
$ws.Range("D32").Value = "'1.99"
$ws.Range("E32").Value = "  +3.48%  "
$ws.Range("D33").Value = "'1.81"
$ws.Range("E33").Value = "  +11.15%  "
$ws.Range("D34").Value = "0.0₃0802"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("D35").Value = "'174.36"
$ws.Range("E35").Value = "  -1.06%  "
$ws.Range("D36").Value = "'4.85"
$ws.Range("E36").Value = "  +9.14%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "'0.403"
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("D39").Value = "'19.03"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").Value = "'1.80"
$ws.Range("E40").Value = "  +6.59%  "
$ws.Range("D41").Value = "'172.58"
$ws.Range("E41").Value = "  +8.91%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "'40.27"
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("D44").Value = "'3.71"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "'22.07"
$ws.Range("E45").Value = "  +4.28%  "
$ws.Range("D46").Value = "'0.0558"
$ws.Range("E46").Value = "  +4.01%  "
$ws.Range("D47").Value = "'0.630"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("D48").Value = "'0.0961"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("E49").Value = "  +1.45%  "
$ws.Range("D50").Value = "'18.52"
$ws.Range("E50").Value = "  +2.42%  "
$ws.Range("D51").Value = "'11.34"
$ws.Range("E51").Value = "  -0.62%  "
